$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.547.96"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "'2.495.59"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'321.94"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'109.09"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'39.32"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "'18.65"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "'2.885.68"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "'2.498.95"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'47.440.76"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'13.36"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("D20").Value = "'6.65"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "'0.0₃0941"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "'2.76"
$ws.Range("E22").Value = "  +15.02%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'246.78"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'25.76"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'9.99"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("D30").Value = "'34.74"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("D32").Value = "'49.93"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'20.30"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'0.0790"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.112"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.46"
$ws.Range("E41").Value = "  +6.84%  "
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "'119.21"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'1.999.41"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").Value = "'2.04"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").Value = "'1.79"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'9.05"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "'56.65"
$ws.Range("E51").Value = "  +3.29%  "
